# Update "想去人数" (want-to-go count) figures for two events that appear
# on both the "展览" sheet and the "全部类型" sheet.
#   F2: 292 -> 294
#   F4: 1120 -> 1134
#   F5/F6: 583 -> 585

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F2").Value = 294
$sheetExhibition.Range("F4").Value = 1134
$sheetExhibition.Range("F5").Value = 585

$sheetAllTypes = $wb.Worksheets.Item("全部类型")
$sheetAllTypes.Range("F2").Value = 294
$sheetAllTypes.Range("F4").Value = 1134
$sheetAllTypes.Range("F6").Value = 585
